$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell with default (unstyled) formatting, used to strip the
# transient Text-number-format style back off D-column cells after
# writing values that would otherwise be auto-coerced to numbers.
$plainStyle = $ws.Range("B2").Style

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.465.97'
$ws.Range('D2').Style = $plainStyle
$ws.Range('E2').Value = '  +2.62%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.677.50'
$ws.Range('D3').Style = $plainStyle
$ws.Range('E3').Value = '  +3.81%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.59'
$ws.Range('D5').Style = $plainStyle
$ws.Range('E5').Value = '  +3.80%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5314'
$ws.Range('D6').Style = $plainStyle
$ws.Range('E6').Value = '  +2.41%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = $plainStyle
$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2674'
$ws.Range('D8').Style = $plainStyle
$ws.Range('E8').Value = '  +4.35%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06395'
$ws.Range('D9').Style = $plainStyle
$ws.Range('E9').Value = '  +2.32%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.59'
$ws.Range('D10').Style = $plainStyle
$ws.Range('E10').Value = '  +6.48%  '

$ws.Range('E11').Value = '  +3.72%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.700.85'
$ws.Range('D12').Style = $plainStyle
$ws.Range('E12').Value = '  +5.08%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.503'
$ws.Range('D13').Style = $plainStyle
$ws.Range('E13').Value = '  +3.41%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5567'
$ws.Range('D14').Style = $plainStyle
$ws.Range('E14').Value = '  +2.69%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0₅8352'
$ws.Range('D15').Style = $plainStyle

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.67'
$ws.Range('D16').Style = $plainStyle
$ws.Range('E16').Value = '  +2.89%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.501.74'
$ws.Range('D17').Style = $plainStyle
$ws.Range('E17').Value = '  +2.82%  '

$ws.Range('E18').Value = '  +0.04%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.773'
$ws.Range('D19').Style = $plainStyle
$ws.Range('E19').Value = '  +3.23%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '195.08'
$ws.Range('D20').Style = $plainStyle
$ws.Range('E20').Value = '  +6.44%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.37'
$ws.Range('D21').Style = $plainStyle
$ws.Range('E21').Value = '  +3.62%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.333'
$ws.Range('D22').Style = $plainStyle
$ws.Range('E22').Value = '  +4.77%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('D23').Style = $plainStyle
$ws.Range('E23').Value = '  +0.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.81'
$ws.Range('D24').Style = $plainStyle
$ws.Range('E24').Value = '  -0.35%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1277'
$ws.Range('D25').Style = $plainStyle
$ws.Range('E25').Value = '  +6.68%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.431'
$ws.Range('D26').Style = $plainStyle
$ws.Range('E26').Value = '  +1.22%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.31'
$ws.Range('D27').Style = $plainStyle
$ws.Range('E27').Value = '  +5.62%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.428'
$ws.Range('D28').Style = $plainStyle
$ws.Range('E28').Value = '  +5.51%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06130'
$ws.Range('D29').Style = $plainStyle
$ws.Range('E29').Value = '  +4.85%  '

$ws.Range('E30').Value = '  +3.15%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.626'
$ws.Range('D31').Style = $plainStyle
$ws.Range('E31').Value = '  +7.77%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.452'
$ws.Range('D32').Style = $plainStyle
$ws.Range('E32').Value = '  +3.55%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.692'
$ws.Range('D33').Style = $plainStyle
$ws.Range('E33').Value = '  +6.19%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.006'
$ws.Range('D34').Style = $plainStyle
$ws.Range('E34').Value = '  +4.20%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.424'
$ws.Range('D35').Style = $plainStyle

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.780'
$ws.Range('D36').Style = $plainStyle
$ws.Range('E36').Value = '  +2.62%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5738'
$ws.Range('D37').Style = $plainStyle
$ws.Range('E37').Value = '  +0.30%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01641'
$ws.Range('D38').Style = $plainStyle
$ws.Range('E38').Value = '  +3.83%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.044'
$ws.Range('D39').Style = $plainStyle
$ws.Range('E39').Value = '  +6.76%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.074.81'
$ws.Range('D40').Style = $plainStyle
$ws.Range('E40').Value = '  +5.96%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8585'
$ws.Range('D41').Style = $plainStyle
$ws.Range('E41').Value = '  +2.31%  '

$ws.Range('E42').Value = '  -0.19%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.10'
$ws.Range('D43').Style = $plainStyle
$ws.Range('E43').Value = '  +0.93%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.825.60'
$ws.Range('D44').Style = $plainStyle
$ws.Range('E44').Value = '  +3.61%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.10'
$ws.Range('D45').Style = $plainStyle
$ws.Range('E45').Value = '  +5.44%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈106'
$ws.Range('D46').Style = $plainStyle
$ws.Range('E46').Value = '  -2.13%  '

$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.003'
$ws.Range('D47').Style = $plainStyle
$ws.Range('E47').Value = '  +0.32%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.120'
$ws.Range('D48').Style = $plainStyle
$ws.Range('E48').Value = '  +3.14%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05206'
$ws.Range('D49').Style = $plainStyle
$ws.Range('E49').Value = '  +1.07%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.479'
$ws.Range('D50').Style = $plainStyle
$ws.Range('E50').Value = '  +8.15%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.033'
$ws.Range('D51').Style = $plainStyle
$ws.Range('E51').Value = '  +4.21%  '
